# adjust aug barrel damage and fix m4 firerates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Aug barrel damage (bullet_deviation col L) and bullet_damage (col I) tweaks ---
$ws.Range("I3").Value = 0.4
$ws.Range("L3").Value = -0.06

$ws.Range("I4").Value = 0.26
$ws.Range("L4").Value = -0.06

$ws.Range("I5").Value = 0.26
$ws.Range("L5").Value = -0.06

$ws.Range("C6").Value = -6
$ws.Range("I6").Value = 0.16
$ws.Range("L6").Value = -0.03

# --- M4 fire rate fixes (col J = fire_rate) ---
$ws.Range("J7").Value = 140
$ws.Range("L7").Value = -0.03

$ws.Range("J8").Value = 100

$ws.Range("J9").Value = 10

$ws.Range("J10").Value = -10

# Restore the active cell selection to where the author ended up
$ws.Range("M18").Select()
